$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 10-11: enter column A values, then column B values
$ws.Cells.Item(10, 1).Value = "我的收藏-店铺"
$ws.Cells.Item(11, 1).Value = "我的收藏-商品"
$ws.Cells.Item(10, 2).Value = "collect-store.html"
$ws.Cells.Item(11, 2).Value = "collect-comm.html"

# Row 12: typed normally
$ws.Cells.Item(12, 1).Value = "我的地址"
$ws.Cells.Item(12, 2).Value = "my-site.html"

# Row 13: typed normally
$ws.Cells.Item(13, 1).Value = "我的地址-收货地址修改"
$ws.Cells.Item(13, 2).Value = "site-change.html"

# Row 14: typed normally
$ws.Cells.Item(14, 1).Value = "我的地址-添加地址"
$ws.Cells.Item(14, 2).Value = "site-add.html"

# Rows 15-16: enter column B values, then column A values
$ws.Cells.Item(15, 2).Value = "bill.html"
$ws.Cells.Item(16, 2).Value = "bill-add.html"
$ws.Cells.Item(15, 1).Value = "发票"
$ws.Cells.Item(16, 1).Value = "新增发票"

# Row 17: typed normally
$ws.Cells.Item(17, 1).Value = "意见反馈"
$ws.Cells.Item(17, 2).Value = "opinion.html"

# Leave the active cell on B17, matching the last edited cell
$null = $ws.Range("B17").Select()
